# Updates cryptos price (D) and 1h volume change (E) columns
# as scraped on Wed Feb 28 16:31:46 UTC 2024 (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.379.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +7.86%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.380.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.12%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "418.30"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "115.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.65%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.596"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.52%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.647"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.105"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.30%  "

$ws.Range("E12").Value = "  +1.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.904.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.09%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.415.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.26%  "

$ws.Range("E17").Value = "  +2.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.135.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000115"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "304.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.91%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.29%  "

$ws.Range("E30").Value = "  +6.72%  "

$ws.Range("E31").Value = "  +6.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +23.59%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "40.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.68%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0511"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.95%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.70%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "137.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.73%  "

$ws.Range("E42").Value = "  +3.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.294"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.51%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.176.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.72%  "

$ws.Range("E50").Value = "  +1.74%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.99"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.54%  "
